# The deck's theme (ppt/theme/theme1.xml, "Integral" / "Red Violet" colours)
# is swapped for the stock "Office Theme" / "Office" colour scheme that
# previously only lived in the Notes Master's theme part (ppt/theme/theme2.xml).
#
# Index order for ThemeColorScheme.Item(n) follows the standard Office
# theme colour order:
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1 6 accent2 7 accent3
#   8 accent4 9 accent5 10 accent6 11 hlink 12 folHlink

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$scheme = $design.SlideMaster.Theme.ThemeColorScheme

$scheme.Item(1).RGB  = 0        # dk1      000000
$scheme.Item(2).RGB  = 16777215 # lt1      FFFFFF
$scheme.Item(3).RGB  = 6968388  # dk2      44546A
$scheme.Item(4).RGB  = 15132391 # lt2      E7E6E6
$scheme.Item(5).RGB  = 13998939 # accent1  5B9BD5
$scheme.Item(6).RGB  = 3243501  # accent2  ED7D31
$scheme.Item(7).RGB  = 10855845 # accent3  A5A5A5
$scheme.Item(8).RGB  = 49407    # accent4  FFC000
$scheme.Item(9).RGB  = 12874308 # accent5  4472C4
$scheme.Item(10).RGB = 4697456  # accent6  70AD47
$scheme.Item(11).RGB = 12673797 # hlink    0563C1
$scheme.Item(12).RGB = 7491477  # folHlink 954F72
